$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.554.82"

$ws.Range("D3").Value = "'2.003.48"
$ws.Range("E3").Value = "'  -4.19%  "

$ws.Range("D4").Value = "'1.015"
$ws.Range("E4").Value = "'  +1.33%  "

$ws.Range("D5").Value = "'329.79"
$ws.Range("E5").Value = "'  -3.84%  "

$ws.Range("E6").Value = "'  +1.19%  "

$ws.Range("D7").Value = "'0.5005"
$ws.Range("E7").Value = "'  -4.39%  "

$ws.Range("D8").Value = "'0.4218"
$ws.Range("E8").Value = "'  -4.62%  "

$ws.Range("D9").Value = "'54.56"
$ws.Range("E9").Value = "'  +0.04%  "

$ws.Range("D10").Value = "'0.09024"
$ws.Range("E10").Value = "'  -3.27%  "

$ws.Range("D11").Value = "'1.117"

$ws.Range("B12").Value = "'WrappedEther"
$ws.Range("C12").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'2.101.87"
$ws.Range("E12").Value = "'  +1.41%  "

$ws.Range("B13").Value = "'Solana"
$ws.Range("C13").Value = "'https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "'23.32"
$ws.Range("E13").Value = "'  -6.18%  "

$ws.Range("D14").Value = "'8.041"
$ws.Range("E14").Value = "'  -6.39%  "

$ws.Range("D15").Value = "'6.465"
$ws.Range("E15").Value = "'  -6.25%  "

$ws.Range("E16").Value = "'  +1.17%  "

$ws.Range("D17").Value = "'94.44"
$ws.Range("E17").Value = "'  -6.69%  "

$ws.Range("E18").Value = "'  -3.85%  "

$ws.Range("D19").Value = "'0.06676"
$ws.Range("E19").Value = "'  +0.14%  "

$ws.Range("D20").Value = "'19.64"
$ws.Range("E20").Value = "'  -7.01%  "

$ws.Range("E21").Value = "'  +1.17%  "

$ws.Range("D22").Value = "'5.974"
$ws.Range("E22").Value = "'  -5.67%  "

$ws.Range("D23").Value = "'29.610.90"

$ws.Range("D24").Value = "'11.99"
$ws.Range("E24").Value = "'  -4.36%  "

$ws.Range("D25").Value = "'2.304"
$ws.Range("E25").Value = "'  +0.13%  "

$ws.Range("D26").Value = "'158.71"
$ws.Range("E26").Value = "'  -2.48%  "

$ws.Range("D27").Value = "'20.73"
$ws.Range("E27").Value = "'  -4.94%  "

$ws.Range("D28").Value = "'6.351"
$ws.Range("E28").Value = "'  -6.29%  "

$ws.Range("D29").Value = "'2.297"
$ws.Range("E29").Value = "'  -8.65%  "

$ws.Range("D30").Value = "'128.15"
$ws.Range("E30").Value = "'  -3.68%  "

$ws.Range("D31").Value = "'1.055"

$ws.Range("D32").Value = "'0.09952"
$ws.Range("E32").Value = "'  -4.78%  "

$ws.Range("D33").Value = "'1.565"
$ws.Range("E33").Value = "'  -6.56%  "

$ws.Range("D34").Value = "'5.835"
$ws.Range("E34").Value = "'  -6.55%  "

$ws.Range("D35").Value = "'3.801"
$ws.Range("E35").Value = "'  -1.57%  "

$ws.Range("D36").Value = "'0.02467"
$ws.Range("E36").Value = "'  -6.30%  "

$ws.Range("D37").Value = "'9.278"
$ws.Range("E37").Value = "'  -9.09%  "

$ws.Range("B38").Value = "'TrustWalletToken"
$ws.Range("C38").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.308"
$ws.Range("E38").Value = "'  -2.71%  "

$ws.Range("B39").Value = "'Hedera"
$ws.Range("C39").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.06401"
$ws.Range("E39").Value = "'  -6.43%  "

$ws.Range("D40").Value = "'0.6549"
$ws.Range("E40").Value = "'  -6.23%  "

$ws.Range("D41").Value = "'11.67"
$ws.Range("E41").Value = "'  -6.71%  "

$ws.Range("D42").Value = "'0.2045"
$ws.Range("E42").Value = "'  -7.38%  "

$ws.Range("E43").Value = "'  +1.23%  "

$ws.Range("D44").Value = "'0.6346"
$ws.Range("E44").Value = "'  -7.01%  "

$ws.Range("D45").Value = "'13.51"
$ws.Range("E45").Value = "'  -6.02%  "

$ws.Range("D46").Value = "'2.193"
$ws.Range("E46").Value = "'  -6.54%  "

$ws.Range("D47").Value = "'1.303"
$ws.Range("E47").Value = "'  -5.04%  "

$ws.Range("D48").Value = "'3.509"
$ws.Range("E48").Value = "'  -3.47%  "

$ws.Range("D49").Value = "'0.00000000341"
$ws.Range("E49").Value = "'  -2.17%  "

$ws.Range("D50").Value = "'0.06987"
$ws.Range("E50").Value = "'  -3.40%  "

$ws.Range("D51").Value = "'1.129"
$ws.Range("E51").Value = "'  -6.74%  "
